# Repair bug edit: add two new data rows to "parsed mile posts" sheet and
# add a new shared comment string ("Intersection with I90").
#
# Net effect on the OOXML:
#   - xl/sharedStrings.xml gains one new unique string, inserted right
#     after "IS" (index 8) -> "Intersection with I90" (index 9). Every
#     shared string used afterwards (the "definitions" sheet) shifts by +1.
#   - xl/worksheets/sheet1.xml ("parsed mile posts") gets two new rows
#     (3 and 4), both cloning row 2's data, with row 3 carrying a comment
#     in column I and row 4 carrying an (empty) formatted cell in column H.
#   - dimension / selected cell on sheet1 move accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parsed mile posts")
$dataRowHeight = $ws.Rows.Item(2).RowHeight

# --- Row 3: clone of row 2 plus a comment in column I ------------------
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 10.15
$ws.Range("C3").Value = 11.56
$ws.Range("D3").Value = 50000
$ws.Range("E3").Value = "IS"
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = "Intersection with I90"
$ws.Rows.Item(3).RowHeight = $dataRowHeight

# --- Row 4: another clone of row 2, with an empty (but styled) H4 ------
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 10.15
$ws.Range("C4").Value = 11.56
$ws.Range("D4").Value = 50000
$ws.Range("E4").Value = "IS"
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 1
# H4 stays empty but picks up the same centered "s=1" format used by A2/A3/A4
$ws.Range("H4").HorizontalAlignment = -4108
$ws.Rows.Item(4).RowHeight = $dataRowHeight

# The committer's cursor ended on F11 when the edit was saved.
$ws.Range("F11").Select()
